# Automatische test-sync: 2025-06-24 19:59:50
# Adds the new "Inlogproblemen" mail-log entry (row 13) to the Logs sheet,
# rolls the corresponding category count into the Dashboard sheet (row 8),
# and extends the conditional formatting + chart series ranges to include
# the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row (row 13)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Inlogproblemen"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Ik kan niet meer inloggen op mijn account. Kunnen jullie helpen?"
$logs.Range("D13").Value = "IT / Technisch probleem"
$logs.Range("E13").Value = "Beste klant,`nBedankt voor je bericht. Om je verder te kunnen helpen met het inlogprobleem, heb ik wat meer informatie nodig. Zou je alsjeblieft je gebruikersnaam willen doorgeven? Hiermee kunnen we het probleem verder onderzoeken en een oplossing bieden.`nAlvast bedankt voor je medewerking.`nMet vriendelijke groet,`n[Naam assistent]`nNederlandse e-mailassistent"
$logs.Range("F13").Value = "2025-06-24 19:59:32"
$logs.Range("G13").Value = "Ja"

# Extend the existing conditional-formatting blocks (D2:D12 -> D2:D13,
# G2:G12 -> G2:G13) in place so the rule set (priorities / dxfIds /
# formulas) stays identical - just the applied range grows.
$dRules = $logs.Range("D2:D12").FormatConditions
for ($i = 1; $i -le $dRules.Count; $i++) {
    $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D13"))
}

$gRules = $logs.Range("G2:G12").FormatConditions
for ($i = 1; $i -le $gRules.Count; $i++) {
    $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G13"))
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: add the new category count (row 8)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "IT / Technisch probleem"
$dash.Range("B8").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category / value series references to A2:A8 / B2:B8
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$8"
$series.Values = "='Dashboard'!`$B`$2:`$B`$8"
